$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) - force text via leading apostrophe then reset style to avoid numeric coercion residue
$ws.Range("D2").Value = "'26.029.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.646.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'215.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.5227"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'0.2609"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06367"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'20.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07681"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.653.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'4.426"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'1.867.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.5543"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8333"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'64.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'26.034.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = "'4.715"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'188.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'6.261"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'144.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1222"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'7.402"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'15.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Value = "'0.05969"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Value = "'3.400"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'3.406"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.651"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.9958"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'2.393"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'2.754"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.5630"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Value = "'5.858"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.8547"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'1.027.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'98.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'1.795.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'55.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Value = "'8.095"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.05144"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.4216"
$ws.Range("D51").Style = "Normal"

# Volume(1h) column (E) - plain text assignment (percent strings are not numeric-coercible)
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  -5.19%  "
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  -6.12%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -7.67%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -0.50%  "
